# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp string (08:30 -> 09:47)
# - Refresh case counts for several countries with newer data
# - Lituania's new total (1908) overtakes Guinea-Bisau's (1902), so the two
#   countries swap places in the (descending, by total-cases) sorted table:
#   row 120 becomes Lituania's updated data, row 121 becomes Guinea-Bisau's
#   former data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Refresh timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 09:47"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 3695302
$ws.Range("C4").Value = 277
$ws.Range("D4").Value = 1679749
$ws.Range("E4").Value = 1874435

# --- Rusia (row 7) ---
$ws.Range("B7").Value = 759203
$ws.Range("C7").Value = 6406
$ws.Range("D7").Value = 539373
$ws.Range("E7").Value = 207707
$ws.Range("G7").Value = 186
$ws.Range("H7").Value = 12123

# --- Armenia (row 52) ---
$ws.Range("B52").Value = 34001
$ws.Range("C52").Value = 442
$ws.Range("D52").Value = 22492
$ws.Range("E52").Value = 10889
$ws.Range("G52").Value = 13
$ws.Range("H52").Value = 620

# --- El Salvador (row 75) ---
$ws.Range("B75").Value = 11207
$ws.Range("C75").Value = 250
$ws.Range("D75").Value = 6422
$ws.Range("E75").Value = 4476

# --- Hungria (row 99) ---
$ws.Range("B99").Value = 4293
$ws.Range("C99").Value = 14
$ws.Range("D99").Value = 3220
$ws.Range("E99").Value = 478

# --- Estonia (row 117) ---
$ws.Range("B117").Value = 2020
$ws.Range("C117").Value = 4
$ws.Range("D117").Value = 1910
$ws.Range("E117").Value = 41

# --- Guinea-Bisau / Lituania swap (rows 120-121) ---
# Row 120 now holds Lituania's refreshed figures.
$ws.Range("A120").Value = "Lituania"
$ws.Range("B120").Value = 1908
$ws.Range("C120").Value = 6
$ws.Range("D120").Value = 1595
$ws.Range("E120").Value = 234
$ws.Range("H120").Value = 79

# Row 121 now holds Guinea-Bisau's (previous row 120) figures.
$ws.Range("A121").Value = "Guinea-Bisau"
$ws.Range("B121").Value = 1902
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 773
$ws.Range("E121").Value = 1103
$ws.Range("H121").Value = 26

# --- Letonia (row 138) ---
$ws.Range("B138").Value = 1185
$ws.Range("C138").Value = 6
$ws.Range("E138").Value = 132

# --- Taiwan (row 161) ---
$ws.Range("B161").Value = 454
$ws.Range("C161").Value = 2
$ws.Range("E161").Value = 7
